$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values must be force-written as text (matching the
# original inline-string cell type) so Excel does not auto-convert them
# to numbers (which would also strip formatting like trailing zeros).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "63.850.76"
Set-TextValue $ws.Range("E2") "  -6.31%  "
# Row 3
Set-TextValue $ws.Range("D3") "3.279.01"
Set-TextValue $ws.Range("E3") "  -8.68%  "
# Row 4
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.16%  "
# Row 5
Set-TextValue $ws.Range("D5") "177.74"
Set-TextValue $ws.Range("E5") "  -13.32%  "
# Row 6
Set-TextValue $ws.Range("D6") "514.06"
Set-TextValue $ws.Range("E6") "  -9.61%  "
# Row 7
Set-TextValue $ws.Range("D7") "0.588"
Set-TextValue $ws.Range("E7") "  -4.04%  "
# Row 8
Set-TextValue $ws.Range("D8") "3.276.90"
Set-TextValue $ws.Range("E8") "  -8.68%  "
# Row 10
Set-TextValue $ws.Range("D10") "0.611"
Set-TextValue $ws.Range("E10") "  -10.19%  "
# Row 11
Set-TextValue $ws.Range("D11") "57.20"
Set-TextValue $ws.Range("E11") "  -5.54%  "
# Row 12
Set-TextValue $ws.Range("D12") "0.130"
Set-TextValue $ws.Range("E12") "  -12.02%  "
# Row 13
Set-TextValue $ws.Range("D13") "0.0000253"
Set-TextValue $ws.Range("E13") "  -9.78%  "
# Row 14
Set-TextValue $ws.Range("D14") "8.98"
Set-TextValue $ws.Range("E14") "  -12.46%  "
# Row 15
Set-TextValue $ws.Range("D15") "3.789.17"
Set-TextValue $ws.Range("E15") "  -8.94%  "
# Row 16
Set-TextValue $ws.Range("D16") "0.119"
Set-TextValue $ws.Range("E16") "  -5.65%  "
# Row 17
Set-TextValue $ws.Range("D17") "3.272.29"
Set-TextValue $ws.Range("E17") "  -9.02%  "
# Row 18
Set-TextValue $ws.Range("D18") "63.510.53"
Set-TextValue $ws.Range("E18") "  -6.47%  "
# Row 19
Set-TextValue $ws.Range("D19") "17.00"
Set-TextValue $ws.Range("E19") "  -10.11%  "
# Row 20
Set-TextValue $ws.Range("D20") "10.72"
Set-TextValue $ws.Range("E20") "  -12.28%  "
# Row 21
Set-TextValue $ws.Range("D21") "0.939"
Set-TextValue $ws.Range("E21") "  -11.50%  "
# Row 22
Set-TextValue $ws.Range("D22") "366.79"
Set-TextValue $ws.Range("E22") "  -8.82%  "
# Row 23
Set-TextValue $ws.Range("D23") "79.38"
Set-TextValue $ws.Range("E23") "  -6.06%  "
# Row 24
Set-TextValue $ws.Range("D24") "3.63"
Set-TextValue $ws.Range("E24") "  -13.23%  "
# Row 25
Set-TextValue $ws.Range("D25") "10.72"
Set-TextValue $ws.Range("E25") "  -14.52%  "
# Row 26
Set-TextValue $ws.Range("D26") "3.78"
Set-TextValue $ws.Range("E26") "  -2.98%  "
# Row 28
Set-TextValue $ws.Range("D28") "2.61"
Set-TextValue $ws.Range("E28") "  -9.43%  "
# Row 29
Set-TextValue $ws.Range("D29") "11.21"
Set-TextValue $ws.Range("E29") "  -9.90%  "
# Row 30
Set-TextValue $ws.Range("D30") "8.22"
Set-TextValue $ws.Range("E30") "  -10.88%  "
# Row 31
Set-TextValue $ws.Range("D31") "28.22"
Set-TextValue $ws.Range("E31") "  -10.47%  "
# Row 32
Set-TextValue $ws.Range("D32") "633.50"
Set-TextValue $ws.Range("E32") "  -5.24%  "
# Row 33
Set-TextValue $ws.Range("D33") "6.60"
Set-TextValue $ws.Range("E33") "  -14.71%  "
# Row 34
Set-TextValue $ws.Range("D34") "10.96"
Set-TextValue $ws.Range("E34") "  -9.21%  "
# Row 35
Set-TextValue $ws.Range("D35") "58.80"
Set-TextValue $ws.Range("E35") "  -7.21%  "
# Row 36
Set-TextValue $ws.Range("D36") "0.102"
Set-TextValue $ws.Range("E36") "  -9.60%  "
# Row 38
Set-TextValue $ws.Range("D38") "35.52"
Set-TextValue $ws.Range("E38") "  -13.73%  "
# Row 39
Set-TextValue $ws.Range("D39") "0.371"
Set-TextValue $ws.Range("E39") "  -9.08%  "
# Row 40
Set-TextValue $ws.Range("D40") "0.997"
Set-TextValue $ws.Range("E40") "  -0.22%  "
# Row 41
Set-TextValue $ws.Range("D41") "0.121"
Set-TextValue $ws.Range("E41") "  -9.13%  "
# Row 42
Set-TextValue $ws.Range("D42") "2.842.36"
Set-TextValue $ws.Range("E42") "  -10.90%  "
# Row 43
Set-TextValue $ws.Range("B43") "PEPE"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D43") "0.0₃0632"
Set-TextValue $ws.Range("E43") "  -16.16%  "
# Row 44
Set-TextValue $ws.Range("B44") "ThetaToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D44") "2.66"
Set-TextValue $ws.Range("E44") "  -18.55%  "
# Row 45
Set-TextValue $ws.Range("B45") "WEMIXToken"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D45") "2.56"
Set-TextValue $ws.Range("E45") "  -8.58%  "
# Row 46
Set-TextValue $ws.Range("B46") "VeChain"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D46") "0.0380"
Set-TextValue $ws.Range("E46") "  -7.23%  "
# Row 47
Set-TextValue $ws.Range("B47") "Fetch.AI"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D47") "2.28"
Set-TextValue $ws.Range("E47") "  -15.13%  "
# Row 48
Set-TextValue $ws.Range("B48") "Stellar"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D48") "0.123"
Set-TextValue $ws.Range("E48") "  -5.85%  "
# Row 49
Set-TextValue $ws.Range("B49") "Monero"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D49") "132.64"
Set-TextValue $ws.Range("E49") "  -4.41%  "
# Row 50
Set-TextValue $ws.Range("B50") "Stacks"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D50") "2.64"
Set-TextValue $ws.Range("E50") "  -1.82%  "
# Row 51
Set-TextValue $ws.Range("B51") "ApeXProtocol"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D51") "2.79"
Set-TextValue $ws.Range("E51") "  -8.75%  "
